$wb = $excel.ActiveWorkbook

# --- Sheet 1: DQ_Report ---
$ws1 = $wb.Worksheets.Item("DQ_Report")

# Insert a new column before column B (ICD_primaerkode), shifting
# ICD_primaerkode -> C, Orpha_Kode -> D, dq_msg -> E
$ws1.Columns.Item(2).Insert()

# New header for the inserted column
$ws1.Range("B1").Value = "Aufnahmenummer"
$ws1.Range("B1").Font.Bold = $true
$ws1.Range("B1").HorizontalAlignment = -4108

# New "Aufnahmenummer" (case id) values for each data row
$ws1.Range("B2").Value = "F_101645"
$ws1.Range("B3").Value = "F_101646"
$ws1.Range("B4").Value = "F_101648"
$ws1.Range("B5").Value = "F_101649"
$ws1.Range("B6").Value = "F_101650"
$ws1.Range("B7").Value = "F_101651"
$ws1.Range("B8").Value = "F_101651"
$ws1.Range("B9").Value = "F_101653"
$ws1.Range("B10").Value = "F_101654"
$ws1.Range("B11").Value = "F_101655"
$ws1.Range("B12").Value = "F_101656"
$ws1.Range("B13").Value = "F_101757"
$ws1.Range("B14").Value = "F_101658"
$ws1.Range("B15").Value = "F_101660"

# --- Sheet 2: Statistik ---
$ws2 = $wb.Worksheets.Item("Statistik")
$ws2.Range("A1").Value = "inst_id"
$ws2.Range("A2").Value = "260123430-Dali"
$ws2.Range("B2").Value = 3.13
$ws2.Range("C2").Value = 96.87
